$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G values (mirrors column F for most rows, row 7 differs)
$ws.Range("G1").Value = -20
$ws.Range("G2").Value = -176
$ws.Range("G3").Value = -30
$ws.Range("G5").Value = 6472
$ws.Range("G6").Value = -2144
$ws.Range("G7").Value = 730
$ws.Range("G9").Value = -1
$ws.Range("G10").Value = -1
$ws.Range("G11").Value = -1

# Row 14 spans change to 1:7 without adding data - handled implicitly by Excel when dimension grows

# Add new summary formulas in column I
$ws.Range("I19").Formula = "=AVERAGE(B1:O1)"
$ws.Range("I20").Formula = "=AVERAGE(A2:O2)"
$ws.Range("I21").Formula = "=AVERAGE(C3:R3)"
$ws.Range("I23").Formula = "=AVERAGE(A5:P5)"
$ws.Range("I24").Formula = "=AVERAGE(A6:R6)"
$ws.Range("I25").Formula = "=AVERAGE(A7:P7)"

# Update selection to match target state
$ws.Range("E31").Select()
